$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeBook")
$ws.Activate()

# A new reporting day arrived: column Z gets the same data as the previous
# day's column (Y), carried forward with style, then the two counties that
# actually reported new cases are bumped by one.
$ws.Range("Y3:Y24").Copy($ws.Range("Z3:Z24"))
$ws.Range("Z3").Value = 116
$ws.Range("Z11").Value = 51

# Extend the running "total" row into the new column, keeping the same look
# as the rest of that row.
$ws.Range("Y25").Copy($ws.Range("Z25"))
$ws.Range("Z25").Formula = "=SUM(Z3:Z24)"

# Reflect the user having scrolled over / selected the next empty block of
# columns while reviewing the sheet.
$ws.Columns("AC:AF").Select()
